$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: add E8, a "Wingdings" styled checkbox-like cell (empty) ---
$ws.Range("E8").Font.Name = "Wingdings"

# --- Row 14: C14 changes from "Tobi" to "Domi" ---
$ws.Range("C14").Value = "Domi"

# --- New row 17, modeled after row 16's formatting (create its strings first so
#     the shared-string table order matches: "Gruppe löschen..." before "ü") ---
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Gruppe löschen -> anlegen Fehler?"
$ws.Range("C17").Value = "Tobi"

# --- New row 18, just A18 carrying the same numeric-centered style, no value ---
$ws.Range("A16").Copy()
$ws.Range("A18").PasteSpecial(-4122)

# --- E14: styled like E8, value "ü" (added after row 17's strings) ---
$ws.Range("E14").Value = "ü"
$ws.Range("E8").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# --- Update selection to match the author's final cursor position ---
$ws.Range("E15").Select() | Out-Null
